# Update "想去人数" (wish-to-go count) figures in column F across the
# workbook's sheets, refreshing the scraped statistics as per the
# "Update gh-pages to output generated at 456a3b4" commit.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 2565
$ws1.Range("F7").Value = 54
$ws1.Range("F8").Value = 50
$ws1.Range("F10").Value = 1522
$ws1.Range("F12").Value = 618
$ws1.Range("F13").Value = 1491
$ws1.Range("F14").Value = 1292
$ws1.Range("F16").Value = 516
$ws1.Range("F17").Value = 3673
$ws1.Range("F20").Value = 744
$ws1.Range("F21").Value = 7
$ws1.Range("F22").Value = 2158
$ws1.Range("F24").Value = 295
$ws1.Range("F27").Value = 1139
$ws1.Range("F29").Value = 46
$ws1.Range("F30").Value = 1001
$ws1.Range("F31").Value = 987

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F18").Value = 241
$ws2.Range("F19").Value = 183
$ws2.Range("F20").Value = 467

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F5").Value = 514
$ws3.Range("F6").Value = 27

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F11").Value = 514
$ws4.Range("F12").Value = 2565
$ws4.Range("F13").Value = 2565
$ws4.Range("F15").Value = 54
$ws4.Range("F16").Value = 50
$ws4.Range("F24").Value = 1522
$ws4.Range("F27").Value = 1491
$ws4.Range("F29").Value = 1292
$ws4.Range("F30").Value = 516
$ws4.Range("F32").Value = 3674
$ws4.Range("F35").Value = 744
$ws4.Range("F36").Value = 2158
$ws4.Range("F38").Value = 295
$ws4.Range("F40").Value = 1139
$ws4.Range("F43").Value = 241
$ws4.Range("F44").Value = 183
$ws4.Range("F45").Value = 467
$ws4.Range("F47").Value = 46
$ws4.Range("F48").Value = 1001
$ws4.Range("F49").Value = 987
